# chengdu_covid19.xlsx data update (Aug. 04) -- append rows 18-21 to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows: date, cdcon22, cdasy22, cdasytocon22, cdcon15, cdasy15, cdasytocon15
$newRows = @(
    @{ Row = 18; Date = 44773; B = 3; C = 4; D = 0; E = 0; F = 0; G = 0 },
    @{ Row = 19; Date = 44774; B = 9; C = 2; D = 4 },
    @{ Row = 20; Date = 44775; B = 1; C = 5; D = 0 },
    @{ Row = 21; Date = 44776; B = 1; C = 0; D = 1 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Copy the date-formatted style from the row above so the new date cell
    # matches the existing column-A date formatting instead of minting a
    # brand new style entry.
    $ws.Range("A" + ($r - 1)).Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    if ($entry.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $entry.E }
    if ($entry.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $entry.F }
    if ($entry.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $entry.G }
}

$excel.CutCopyMode = $false

# Reflect the scrolled/selected view from the saved workbook.
$ws.Range("E23").Select() | Out-Null
